$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 314 (everything currently at/after row 314
# shifts down by two rows, which also naturally reproduces the extra
# duplicated pair that ends up at the new rows 338-339).
$ws.Rows.Item(314).Insert()
$ws.Rows.Item(314).Insert()

# Fill in the new row 314 (Primera) with the new weekly record.
$ws.Cells.Item(314,1).Value = 8
$ws.Cells.Item(314,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(314,3).Value = "Coquimbo"
$ws.Cells.Item(314,4).Value = 44610
$ws.Cells.Item(314,5).Value = 4
$ws.Cells.Item(314,6).Value = 100112009
$ws.Cells.Item(314,7).Value = "Acelga"
$ws.Cells.Item(314,8).Value = "Sin especificar"
$ws.Cells.Item(314,9).Value = "Primera"
$ws.Cells.Item(314,10).Value = 2500
$ws.Cells.Item(314,11).Value = 500
$ws.Cells.Item(314,12).Value = 600
$ws.Cells.Item(314,13).Value = 550
$ws.Cells.Item(314,14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(314,15).Value = "Provincia del Elquí"
$ws.Cells.Item(314,16).Value = 275
$ws.Cells.Item(314,17).Value = 2
$ws.Cells.Item(314,18).Value = "Hortaliza"

# Fill in the new row 315 (Segunda) with the new weekly record.
$ws.Cells.Item(315,1).Value = 8
$ws.Cells.Item(315,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(315,3).Value = "Coquimbo"
$ws.Cells.Item(315,4).Value = 44610
$ws.Cells.Item(315,5).Value = 4
$ws.Cells.Item(315,6).Value = 100112009
$ws.Cells.Item(315,7).Value = "Acelga"
$ws.Cells.Item(315,8).Value = "Sin especificar"
$ws.Cells.Item(315,9).Value = "Segunda"
$ws.Cells.Item(315,10).Value = 1360
$ws.Cells.Item(315,11).Value = 400
$ws.Cells.Item(315,12).Value = 450
$ws.Cells.Item(315,13).Value = 425
$ws.Cells.Item(315,14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(315,15).Value = "Provincia del Elquí"
$ws.Cells.Item(315,16).Value = 212
$ws.Cells.Item(315,17).Value = 2
$ws.Cells.Item(315,18).Value = "Hortaliza"
